# Weekly fruit/vegetable update: prepend a new week's pair of records
# (Primera / Segunda quality) for Acelga at "Terminal La Palmera de La
# Serena", ahead of the existing data block that starts at row 589.
# Inserting 2 rows there shifts all the old rows (formerly 589-628) down
# to 591-630, growing the used range from A1:R628 to A1:R630.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A589:A590").EntireRow.Insert()

# New row 589 - "Primera" quality
$ws.Cells.Item(589, 1).Value = 8
$ws.Cells.Item(589, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(589, 3).Value = "Coquimbo"
$ws.Cells.Item(589, 4).Value = 45021
$ws.Cells.Item(589, 5).Value = 4
$ws.Cells.Item(589, 6).Value = 100112009
$ws.Cells.Item(589, 7).Value = "Acelga"
$ws.Cells.Item(589, 8).Value = "Sin especificar"
$ws.Cells.Item(589, 9).Value = "Primera"
$ws.Cells.Item(589, 10).Value = 2060
$ws.Cells.Item(589, 11).Value = 500
$ws.Cells.Item(589, 12).Value = 600
$ws.Cells.Item(589, 13).Value = 550
$ws.Cells.Item(589, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(589, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(589, 16).Value = 275
$ws.Cells.Item(589, 17).Value = 2
$ws.Cells.Item(589, 18).Value = "Hortaliza"

# New row 590 - "Segunda" quality
$ws.Cells.Item(590, 1).Value = 8
$ws.Cells.Item(590, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(590, 3).Value = "Coquimbo"
$ws.Cells.Item(590, 4).Value = 45021
$ws.Cells.Item(590, 5).Value = 4
$ws.Cells.Item(590, 6).Value = 100112009
$ws.Cells.Item(590, 7).Value = "Acelga"
$ws.Cells.Item(590, 8).Value = "Sin especificar"
$ws.Cells.Item(590, 9).Value = "Segunda"
$ws.Cells.Item(590, 10).Value = 1500
$ws.Cells.Item(590, 11).Value = 400
$ws.Cells.Item(590, 12).Value = 450
$ws.Cells.Item(590, 13).Value = 425
$ws.Cells.Item(590, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(590, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(590, 16).Value = 212
$ws.Cells.Item(590, 17).Value = 2
$ws.Cells.Item(590, 18).Value = "Hortaliza"
